# "new branch for gulp" -- add createdAt/updatedAt columns to the "Task"
# table header row (row 7) and update the view state (zoom + selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the Task header row (row 7) with two more field-name columns,
# reusing the same shared strings already used by row 16 ("(createdAt)" /
# "(updatedat)"). These new cells pick up the sheet's default (unstyled)
# formatting, matching the rest of the workbook's un-styled cells.
$ws.Range("H7").Value = "(createdAt)"
$ws.Range("I7").Value = "(updatedat)"

# Update the view: zoom out from 200% to 125%, and move the selection /
# scroll position down to B16 (the last header in the sheet).
$win = $excel.ActiveWindow
$win.Zoom = 125
$ws.Range("B16").Select() | Out-Null
